$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 183; existing rows 183-197 shift down to 184-198.
$ws.Rows(183).Insert()

# Populate the newly inserted row 183 with the new weekly data point.
$ws.Cells.Item(183, 1).Value = 5
$ws.Cells.Item(183, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(183, 3).Value = "Maule"
$ws.Cells.Item(183, 4).Value = 44918
$ws.Cells.Item(183, 5).Value = 7
$ws.Cells.Item(183, 6).Value = 100112031
$ws.Cells.Item(183, 7).Value = "Poroto verde"
$ws.Cells.Item(183, 8).Value = "Sin especificar"
$ws.Cells.Item(183, 9).Value = "Primera"
$ws.Cells.Item(183, 10).Value = 300
$ws.Cells.Item(183, 11).Value = 20000
$ws.Cells.Item(183, 12).Value = 20000
$ws.Cells.Item(183, 13).Value = 20000
$ws.Cells.Item(183, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(183, 15).Value = "Región del Maule"
$ws.Cells.Item(183, 16).Value = 800
$ws.Cells.Item(183, 17).Value = 25
$ws.Cells.Item(183, 18).Value = "Hortaliza"
